# Finalização do estudo básico de pandas
# Rename the second worksheet from "Planilha2" to "Resumo".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha2")
$ws.Name = "Resumo"
